# Weekly refresh: insert a new day's price observation as the new first
# record of the "Orégano" series (new row 322) and shift the existing
# historical rows (old 322..389) down by one (new 323..390).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 322; Excel pushes rows 322..389 down to 323..390
# and copies formatting (incl. the date-number style on column D) from
# the row that was previously at this position.
$ws.Rows.Item(322).Insert()

# Populate the newly inserted row with the latest observation.
$ws.Range("A322").Value = 6
$ws.Range("B322").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C322").Value = "Metropolitana"
$ws.Range("D322").Value = 45275
$ws.Range("E322").Value = 13
$ws.Range("F322").Value = 100112029
$ws.Range("G322").Value = "Orégano"
$ws.Range("H322").Value = "Sin especificar"
$ws.Range("I322").Value = "Primera"
$ws.Range("J322").Value = 32
$ws.Range("K322").Value = 16000
$ws.Range("L322").Value = 16000
$ws.Range("M322").Value = 16000
$ws.Range("N322").Value = "$/docena de atados"
$ws.Range("O322").Value = "Región Metropolitana"
$ws.Range("P322").Value = 5333
$ws.Range("Q322").Value = 3
$ws.Range("R322").Value = "Hortaliza"
